$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '72.225.24'
$ws.Range("E2").Value = '  +0.92%  '
$ws.Range("D3").Value = '4.021.20'
$ws.Range("E3").Value = '  -0.16%  '
$ws.Range("E4").Value = '  -0.07%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '527.38'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.99%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '149.83'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +2.08%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.700'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +13.23%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.748'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +2.39%  '
$ws.Range("E10").Value = '  -0.71%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '50.68'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +9.16%  '
$ws.Range("E12").Value = '  -2.44%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '10.74'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.02%  '
$ws.Range("D14").Value = '4.653.90'
$ws.Range("E14").Value = '  -0.62%  '
$ws.Range("D15").Value = '4.009.11'
$ws.Range("E15").Value = '  -0.56%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '14.00'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -1.12%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '20.61'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -2.46%  '
$ws.Range("E18").Value = '  -0.23%  '
$ws.Range("E19").Value = '  -1.37%  '
$ws.Range("D20").Value = '72.098.74'
$ws.Range("E20").Value = '  +0.69%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '429.95'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.15%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '97.49'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +2.63%  '
$ws.Range("E23").Value = '  +0.28%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '4.18'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +3.63%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '14.29'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.30%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '11.12'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -7.23%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.75'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -3.76%  '
$ws.Range("B28").Value = 'Toncoin'
$ws.Range("C28").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '3.71'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +21.56%  '
$ws.Range("B29").Value = 'LEO'
$ws.Range("C29").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.86'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.57%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '36.67'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.13%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.39'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +4.78%  '
$ws.Range("E32").Value = '  +1.99%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '13.41'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.30%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '680.36'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -3.13%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '47.64'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +17.93%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '65.39'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -3.28%  '
$ws.Range("E37").Value = '  +1.56%  '
$ws.Range("E38").Value = '  -1.54%  '
$ws.Range("D39").Value = '0.0₃0823'
$ws.Range("E39").Value = '  -7.91%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.40'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +8.38%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.39'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -9.27%  '
$ws.Range("E42").Value = '  +0.15%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.997'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.30%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0489'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +1.58%  '
$ws.Range("B45").Value = 'Stellar'
$ws.Range("C45").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.149'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +3.96%  '
$ws.Range("B46").Value = 'THORChain'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '10.00'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +11.33%  '
$ws.Range("B47").Value = 'ApeXProtocol'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.43'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -2.70%  '
$ws.Range("B48").Value = 'Fetch.AI'
$ws.Range("C48").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.64'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -4.49%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '3.01'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -4.31%  '
$ws.Range("B50").Value = 'Monero'
$ws.Range("C50").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '144.95'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +2.09%  '
$ws.Range("B51").Value = 'FLOKI'
$ws.Range("C51").Value = 'https://coinranking.com/coin/fmHk13Rqw+floki-floki'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.000267'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -3.16%  '
